# "Fruta / hortaliza, semanal"
# A new weekly price record is inserted as row 11 (pushing the existing
# rows 11-93 down to 12-94). Populate the newly inserted row with the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; Excel shifts rows 11..93 down to 12..94
# and the worksheet dimension grows from R93 to R94 automatically.
$ws.Rows("11").Insert()

$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11, 3).Value = "Ñuble"
$ws.Cells.Item(11, 4).Value = 45111
$ws.Cells.Item(11, 5).Value = 16
$ws.Cells.Item(11, 6).Value = 100112001
$ws.Cells.Item(11, 7).Value = "Berenjena"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 9000
$ws.Cells.Item(11, 12).Value = 9000
$ws.Cells.Item(11, 13).Value = 9000
$ws.Cells.Item(11, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 150
$ws.Cells.Item(11, 17).Value = 60
$ws.Cells.Item(11, 18).Value = "Hortaliza"
